# Append the "Monday, Jan 16" departures (rows 594-613) to the KRK
# Departures "Main Data" sheet, following the existing table layout:
# NUMBER | DATE | TIME | FLIGHT | TO | SHORT | AIRLINE | MODEL |
# AIRCFAT ID | STATUS | DIFFERENCE | (L) | (M)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(593, "Monday, Jan 16", "5:40 AM", "LO159",   "Tel Aviv",    "(TLV)", "LOT ",               "B38M", "(SP-LVC)", "6:13 AM", $null, "0 hours, 33 minutes", $null)
    ,@(594, "Monday, Jan 16", "5:40 AM", "LO3910",  "Warsaw",      "(WAW)", "LOT ",               "E195", "(SP-LNM)", "6:00 AM", $null, "0 hours, 20 minutes", $null)
    ,@(595, "Monday, Jan 16", "5:50 AM", "W65013",  "Leeds",       "(LBA)", "Wizz Air ",          "A21N", "(HA-LZI)", "6:21 AM", $null, "0 hours, 31 minutes", $null)
    ,@(596, "Monday, Jan 16", "6:00 AM", "FR4023",  "Lisbon",      "(LIS)", "Buzz ",              "B38M", "(SP-RZB)", "6:19 AM", $null, "0 hours, 19 minutes", $null)
    ,@(597, "Monday, Jan 16", "6:00 AM", "KL1992",  "Amsterdam",   "(AMS)", "KLM ",               "E295", "(PH-NXN)", "6:15 AM", $null, "0 hours, 15 minutes", $null)
    ,@(598, "Monday, Jan 16", "6:15 AM", "W65033",  "Oslo",        "(OSL)", "Wizz Air ",          "A21N", "(HA-LVG)", "6:42 AM", $null, "0 hours, 27 minutes", $null)
    ,@(599, "Monday, Jan 16", "6:20 AM", "FR6208",  "Madrid",      "(MAD)", "Buzz ",              "B38M", "(SP-RZC)", "6:37 AM", $null, "0 hours, 17 minutes", $null)
    ,@(600, "Monday, Jan 16", "6:20 AM", "LH1371",  "Frankfurt",   "(FRA)", "Lufthansa ",         "A21N", "(D-AIEK)", "6:45 AM", $null, "0 hours, 25 minutes", $null)
    ,@(601, "Monday, Jan 16", "6:25 AM", "FR3850",  "Venice",      "(TSF)", "Buzz ",              "B38M", "(SP-RZF)", "6:50 AM", $null, "0 hours, 25 minutes", $null)
    ,@(602, "Monday, Jan 16", "6:45 AM", "FR2229",  "Valencia",    "(VLC)", "Buzz ",              "B38M", "(SP-RZD)", "7:07 AM", $null, "0 hours, 22 minutes", $null)
    ,@(603, "Monday, Jan 16", "6:45 AM", "FR8672",  "Girona",      "(GRO)", "Ryanair ",           "B738", "(SP-RSA)", "7:15 AM", $null, "0 hours, 30 minutes", $null)
    ,@(604, "Monday, Jan 16", "6:45 AM", "W65047",  "Barcelona",   "(BCN)", "Wizz Air ",          "A21N", "(HA-LVO)", "7:02 AM", $null, "0 hours, 17 minutes", $null)
    ,@(605, "Monday, Jan 16", "6:50 AM", "FR6310",  "Oslo",        "(TRF)", "Ryanair ",           "B38M", "(SP-RZK)", "7:24 AM", $null, "0 hours, 34 minutes", $null)
    ,@(606, "Monday, Jan 16", "7:15 AM", "FR6222",  "Bournemouth", "(BOH)", "Ryanair ",           "B738", "(SP-RKB)", "7:44 AM", $null, "0 hours, 29 minutes", $null)
    ,@(607, "Monday, Jan 16", "7:25 AM", "OS600",   "Vienna",      "(VIE)", "Austrian Airlines ", "E195", "(OE-LWA)", "7:28 AM", $null, "0 hours, 3 minutes",  $null)
    ,@(608, "Monday, Jan 16", "7:30 AM", "LH1625",  "Munich",      "(MUC)", "Lufthansa ",         "A320", "(D-AIWA)", "7:48 AM", $null, "0 hours, 18 minutes", $null)
    ,@(609, "Monday, Jan 16", "8:25 AM", "W65075",  "Eindhoven",   "(EIN)", "Wizz Air ",          "A321", "(HA-LXO)", "8:28 AM", $null, "0 hours, 3 minutes",  $null)
    ,@(610, "Monday, Jan 16", "8:50 AM", "LO3904",  "Warsaw",      "(WAW)", "LOT ",               "E75S", "(SP-LIA)", "9:02 AM", $null, "0 hours, 12 minutes", $null)
    ,@(611, "Monday, Jan 16", "9:10 AM", "FR5626",  "Edinburgh",   "(EDI)", "Buzz ",              "B38M", "(SP-RZH)", "9:41 AM", $null, "0 hours, 31 minutes", $null)
    ,@(612, "Monday, Jan 16", "9:45 AM", "FR1812",  "London",      "(LTN)", "Ryanair ",           "B738", "(EI-EMR)", "9:48 AM", $null, "0 hours, 3 minutes",  $null)
)

$startRow = 594
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 1; $c -le $row.Count; $c++) {
        $val = $row[$c - 1]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}
